$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: the paragraph ending in "...le faltaba." currently carries
# the (hidden) _GoBack bookmark right after its final "." run. Insert
# a new trailing run of spaces (sz=32 half-points = 16pt) there, then
# strip the bookmark from this spot - it gets re-created later, right
# before "BIBLIOGRAFIA".
# ------------------------------------------------------------------
$bmOld = $d.Bookmarks.Item("_GoBack")
$insPoint = $bmOld.Start

$spaces = "                               "
$insRange = $d.Range($insPoint, $insPoint)
$insRange.InsertAfter($spaces)

$newRunRange = $d.Range($insPoint, $insPoint + $spaces.Length)
$newRunRange.Font.Size = 16

$bmOld = $d.Bookmarks.Item("_GoBack")
$bmOld.Delete()

# ------------------------------------------------------------------
# Step 2: locate the "BIBLIOGRAFIA" heading and trim the run of 13
# blank spacer paragraphs in front of it down to just 6 (remove 7).
# ------------------------------------------------------------------
$bibIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "BIBLIOGRAFIA`r") {
        $bibIndex = $i
        break
    }
}

$firstBlankIndex = $bibIndex - 13
for ($i = 0; $i -lt 7; $i++) {
    $p = $d.Paragraphs.Item($firstBlankIndex)
    $p.Range.Delete()
}

# ------------------------------------------------------------------
# Step 3: re-create the _GoBack bookmark, collapsed, immediately
# before the "BIBLIOGRAFIA" run (start of its paragraph).
# ------------------------------------------------------------------
$bibIndex2 = -1
$count2 = $d.Paragraphs.Count
for ($i = 1; $i -le $count2; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "BIBLIOGRAFIA`r") {
        $bibIndex2 = $i
        break
    }
}

$pBib = $d.Paragraphs.Item($bibIndex2)
$startPos = $pBib.Range.Start
$rBookmark = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $rBookmark)
